# Einstein_Turtlebot_Presentaion.pptx - "Add files via upload"
#
# The "FUTURE PLANS" slide (slide 8) has a bulleted body placeholder.
# The bullet that read:
#   "Einstein raise arms and laugh based on vision and commands"
# is changed to:
#   "Einstein raise arms and dance based on vision and commands"
# (the single word "laugh" -> "dance"), which PowerPoint represents as
# the paragraph being split into three runs because the edited word sits
# in the middle of the original run:
#   "Einstein raise arms " / "and dance based " / "on vision and commands"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$paragraphs = $tr.Paragraphs()

$target = $null
foreach ($para in $paragraphs) {
    if ($para.Text -like "Einstein raise arms and laugh based on vision and commands*") {
        $target = $para
    }
}

if ($target -ne $null) {
    # "Einstein raise arms " = chars 1-20, "and laugh based " = chars 21-36
    $middle = $target.Characters(21, 16)
    $middle.Text = "and dance based "
}
